$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 75 - this shifts existing rows 75..168 down to 76..169
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new record's data
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = 44483
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100112008
$ws.Range("G75").Value = "Coliflor"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 300
$ws.Range("K75").Value = 750
$ws.Range("L75").Value = 800
$ws.Range("M75").Value = 775
$ws.Range("N75").Value = "`$/unidad"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 775
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"
